$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1819.9231
$ws.Range("I38").Value = 1891.6666
$ws.Range("J38").Value = 959
$ws.Range("K38").Value = 5674.9998
$ws.Range("L38").Value = 2877
$ws.Range("M38").Value = -5302.9998
$ws.Range("N38").Value = -3621

$ws.Range("H61").Value = 2876.5557
$ws.Range("I61").Value = 296.33334
$ws.Range("J61").Value = 4166.6665
$ws.Range("K61").Value = 889.0000200000001
$ws.Range("L61").Value = 12499.9995
$ws.Range("M61").Value = -717.0000200000001
$ws.Range("N61").Value = -12843.9995

$ws.Range("H103").Value = 816.8570999999999
$ws.Range("I103").Value = 849.6111
$ws.Range("K103").Value = 2548.8333
$ws.Range("M103").Value = -1962.8333

$ws.Range("H132").Value = 13923.045
$ws.Range("I132").Value = 2339.7222
$ws.Range("J132").Value = 60256.332
$ws.Range("K132").Value = 7019.1666
$ws.Range("L132").Value = 180768.996
$ws.Range("M132").Value = -4489.1666
$ws.Range("N132").Value = -185828.996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3260.5
$ws.Range("I45").Value = 3136.6
$ws.Range("K45").Value = 3136.6
$ws.Range("M45").Value = -2759.6

$ws.Range("H61").Value = 3628.258
$ws.Range("I61").Value = 2728.1667
$ws.Range("K61").Value = 2728.1667
$ws.Range("M61").Value = -2516.1667

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 3901.3125
$ws.Range("I122").Value = 3176.258
$ws.Range("K122").Value = 9528.773999999999
$ws.Range("M122").Value = -7078.773999999999

$ws.Range("H132").Value = 2780.1428
$ws.Range("I132").Value = 2694.4
$ws.Range("K132").Value = 8083.200000000001
$ws.Range("M132").Value = -5553.200000000001

$ws.Range("H136").Value = 3628.258
$ws.Range("I136").Value = 2728.1667
$ws.Range("K136").Value = 8184.500100000001
$ws.Range("M136").Value = -5634.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 64991.812
$ws.Range("J20").Value = 102247.1
$ws.Range("L20").Value = 102247.1
$ws.Range("N20").Value = -102741.1

$ws.Range("H99").Value = 23232.363
$ws.Range("I99").Value = 31265.143
$ws.Range("J99").Value = 9175
$ws.Range("K99").Value = 31265.143
$ws.Range("L99").Value = 9175
$ws.Range("M99").Value = -29767.143
$ws.Range("N99").Value = -12171

$ws.Range("H105").Value = 4328.4443
$ws.Range("I105").Value = 4328.4443
$ws.Range("K105").Value = 4328.4443
$ws.Range("M105").Value = -2581.4443

$ws.Range("H109").Value = 41666.668
$ws.Range("J109").Value = 41666.668
$ws.Range("L109").Value = 41666.668
$ws.Range("N109").Value = -44440.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2109.9443
$ws.Range("I16").Value = 828.2308
$ws.Range("J16").Value = 5442.4
$ws.Range("K16").Value = 828.2308
$ws.Range("L16").Value = 5442.4
$ws.Range("M16").Value = -541.2308
$ws.Range("N16").Value = -6016.4

$ws.Range("H31").Value = 4103
$ws.Range("I31").Value = 4220.0967
$ws.Range("J31").Value = 3699.6667
$ws.Range("K31").Value = 4220.0967
$ws.Range("L31").Value = 3699.6667
$ws.Range("M31").Value = -3925.0967
$ws.Range("N31").Value = -4289.6667

$ws.Range("H34").Value = 4103
$ws.Range("I34").Value = 4220.0967
$ws.Range("J34").Value = 3699.6667
$ws.Range("K34").Value = 4220.0967
$ws.Range("L34").Value = 3699.6667
$ws.Range("M34").Value = -4018.0967
$ws.Range("N34").Value = -4103.6667

$ws.Range("H58").Value = 3123.0908
$ws.Range("I58").Value = 2723.9167
$ws.Range("J58").Value = 3602.1
$ws.Range("K58").Value = 2723.9167
$ws.Range("L58").Value = 3602.1
$ws.Range("M58").Value = -2520.9167
$ws.Range("N58").Value = -4008.1

$ws.Range("H113").Value = 2109.9443
$ws.Range("I113").Value = 828.2308
$ws.Range("J113").Value = 5442.4
$ws.Range("K113").Value = 828.2308
$ws.Range("L113").Value = 5442.4
$ws.Range("M113").Value = 1341.7692
$ws.Range("N113").Value = -9782.4

$ws.Range("H132").Value = 18063.695
$ws.Range("I132").Value = 23004.8
$ws.Range("J132").Value = 8799.125
$ws.Range("K132").Value = 69014.39999999999
$ws.Range("L132").Value = 26397.375
$ws.Range("M132").Value = -66484.39999999999
$ws.Range("N132").Value = -31457.375

$ws.Range("H134").Value = 4085.4285
$ws.Range("I134").Value = 3988.5557
$ws.Range("K134").Value = 11965.6671
$ws.Range("M134").Value = -9430.667099999999

$ws.Range("H135").Value = 59999.77
$ws.Range("J135").Value = 59999.77
$ws.Range("L135").Value = 59999.77
$ws.Range("N135").Value = -70139.76999999999

$ws.Range("H136").Value = 3123.0908
$ws.Range("I136").Value = 2723.9167
$ws.Range("J136").Value = 3602.1
$ws.Range("K136").Value = 8171.750100000001
$ws.Range("L136").Value = 10806.3
$ws.Range("M136").Value = -5621.750100000001
$ws.Range("N136").Value = -15906.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 4999.5
$ws.Range("I100").Value = 4999.5
$ws.Range("K100").Value = 14998.5
$ws.Range("M100").Value = -14187.5

$ws.Range("H107").Value = 1543.1613
$ws.Range("J107").Value = 1671.6666
$ws.Range("L107").Value = 5014.9998
$ws.Range("N107").Value = -8854.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 222.875
$ws.Range("I2").Value = 91.666664
$ws.Range("K2").Value = 91.666664
$ws.Range("M2").Value = 21.333336

$ws.Range("H70").Value = 90811.92
$ws.Range("I70").Value = 143843.25
$ws.Range("J70").Value = 5961.8
$ws.Range("K70").Value = 143843.25
$ws.Range("L70").Value = 5961.8
$ws.Range("M70").Value = -143573.25
$ws.Range("N70").Value = -6501.8

$ws.Range("H73").Value = 90811.92
$ws.Range("I73").Value = 143843.25
$ws.Range("J73").Value = 5961.8
$ws.Range("K73").Value = 143843.25
$ws.Range("L73").Value = 5961.8
$ws.Range("M73").Value = -142907.25
$ws.Range("N73").Value = -7833.8

$ws.Range("H121").Value = 50000000
$ws.Range("J121").Value = 50000000
$ws.Range("L121").Value = 50000000
$ws.Range("N121").Value = -50003494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2514.3845
$ws.Range("I68").Value = 2482.3333
$ws.Range("K68").Value = 2482.3333
$ws.Range("M68").Value = -1733.3333

$ws.Range("H71").Value = 2514.3845
$ws.Range("I71").Value = 2482.3333
$ws.Range("K71").Value = 12411.6665
$ws.Range("M71").Value = -8667.666499999999

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4169.579
$ws.Range("I132").Value = 3988.2
$ws.Range("J132").Value = 4849.75
$ws.Range("K132").Value = 11964.6
$ws.Range("L132").Value = 14549.25
$ws.Range("M132").Value = -9434.599999999999
$ws.Range("N132").Value = -19609.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 655.8
$ws.Range("I100").Value = 726.3333
$ws.Range("J100").Value = 550
$ws.Range("K100").Value = 1452.6666
$ws.Range("L100").Value = 1100
$ws.Range("M100").Value = -911.6666
$ws.Range("N100").Value = -2182

$ws.Range("H122").Value = 1615.8667
$ws.Range("I122").Value = 1374.1
$ws.Range("J122").Value = 2099.4
$ws.Range("K122").Value = 4122.299999999999
$ws.Range("L122").Value = 6298.200000000001
$ws.Range("M122").Value = -1672.299999999999
$ws.Range("N122").Value = -11198.2

$ws.Range("H136").Value = 3651.5938
$ws.Range("I136").Value = 3586.8262
$ws.Range("J136").Value = 3817.111
$ws.Range("K136").Value = 10760.4786
$ws.Range("L136").Value = 11451.333
$ws.Range("M136").Value = -8210.4786
$ws.Range("N136").Value = -16551.333

$ws.Range("H138").Value = 91666.664
$ws.Range("J138").Value = 91666.664
$ws.Range("L138").Value = 91666.664
$ws.Range("N138").Value = -101946.664
